$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Förändrad" (Changed) date column C was bumped from 45171 (2023-09-02)
# to 45172 (2023-09-03) for every data row (rows 2 through 348).
$ws.Range("C2:C348").Value = 45172
